$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10, shifting existing rows 10-19 down to 11-20.
$ws.Rows("10:10").Insert()

# Fill in the new row 10 with the node data (same text used elsewhere in the
# logic tree, so it reuses existing shared strings).
$ws.Range("A10").Value = "Problem:Is engine misfire on multiple cylinders?(OBD Codes indicating multiple cylinders) (Please answer as: Yes)"
$ws.Range("B10").Value = "Possible_Problem"
$ws.Range("C10").Value = "Possible_Problem:40% Ignition Coil`n30% Fuel Injector`n15% Excessive Carbon Buildup (Damaged Valve)`n10% Spark Plug`n5% VVT Actuator"

# Match the formatting used by the other "Possible_Problem" answer cells in
# column C (wrapped text).
$ws.Range("C10").WrapText = $true

# Row height for the new row (matches the equivalent content row elsewhere).
$ws.Rows("10:10").RowHeight = 259.2

# Update the view so the newly edited area is visible/selected.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C10").Select()
